# "Switching back to LDA for a second"
# Adds a third worksheet ("New Title") with an LDA classification-report
# table (mirroring the layout of the "Titles" sheet), selects it as the
# active sheet, and updates the "Titles" sheet selection.

$wb = $excel.ActiveWorkbook

$doc    = $wb.Worksheets.Item("Document")
$titles = $wb.Worksheets.Item("Titles")

# ---------------------------------------------------------------------
# 1. Update the selection on "Titles" before we add the new sheet.
# ---------------------------------------------------------------------
$titles.Range("F1:H11").Select()

# ---------------------------------------------------------------------
# 2. Add the new worksheet after "Titles" (becomes the 3rd / last sheet
#    and is made the active sheet, matching activeTab going 1 -> 2).
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $titles)
$newSheet.Name = "New Title"

# ---------------------------------------------------------------------
# 3. Populate the header row.
# ---------------------------------------------------------------------
$newSheet.Range("F1").Value = "effect"
$newSheet.Range("G1").Value = "weighted "
$newSheet.Range("H1").Value = "Feasible increase"

# ---------------------------------------------------------------------
# 4. Populate the per-class rows (A..E) - LDA classification report.
# ---------------------------------------------------------------------
$rows = @(
    @{ r=2;  label="earn";              B=0.97; C=0.99; D=0.98;                 E=1077 },
    @{ r=3;  label="money-fx";          B=0.67; C=0.59; D=0.63;                 E=87   },
    @{ r=4;  label="money-supply";      B=0.79; C=0.82; D=0.81;                 E=28   },
    @{ r=5;  label="trade";             B=0.87; C=1;    D=0.93;                 E=76   },
    @{ r=6;  label="acq";               B=0.95; C=0.96; D=0.96;                 E=695  },
    @{ r=7;  label="grain-wheat";       B=1;    C=1;    D=1;                    E=35   },
    @{ r=8;  label="interest";          B=0.78; C=0.76; D=0.77;                 E=82   },
    @{ r=9;  label="crude";             B=0.93; C=0.86; D=0.89;                 E=119  },
    @{ r=10; label="interest-money-fx"; B=0.34; C=0.25; D=0.28999999999999998;  E=40   },
    @{ r=11; label="ship";              B=0.96; C=0.69; D=0.81;                 E=36   }
)

foreach ($row in $rows) {
    $r = $row.r
    $newSheet.Range("A$r").Value = $row.label
    $newSheet.Range("B$r").Value = $row.B
    $newSheet.Range("C$r").Value = $row.C
    $newSheet.Range("D$r").Value = $row.D
    $newSheet.Range("E$r").Value = $row.E
    $newSheet.Range("F$r").Formula = "=E$r/SUM(E`$2:E`$11)"
    $newSheet.Range("G$r").Formula = "=D$r*F$r"
    $newSheet.Range("H$r").Formula = "=(1-D$r)*F$r"
}

# Column widths for the percentage columns (bestFit in the source file).
$newSheet.Columns.Item(3).ColumnWidth = 5.1640625
$newSheet.Columns.Item(4).ColumnWidth = 5.1640625
$newSheet.Columns.Item(5).ColumnWidth = 5.1640625

# ---------------------------------------------------------------------
# 5. Selection on the new sheet.
# ---------------------------------------------------------------------
$newSheet.Range("H3").Select()
